$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(32).Insert()

$ws.Range("A32").Value = 6
$ws.Range("B32").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C32").Value = "Metropolitana"
$ws.Range("D32").Value = 45099
$ws.Range("E32").Value = 13
$ws.Range("F32").Value = 100114007
$ws.Range("G32").Value = "Jengibre"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 190
$ws.Range("K32").Value = 20000
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = 20000
$ws.Range("N32").Value = "`$/caja 13 kilos"
$ws.Range("O32").Value = "Perú"
$ws.Range("P32").Value = 1538
$ws.Range("Q32").Value = 13
$ws.Range("R32").Value = "Hortaliza"
